$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.419.84'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.377.62'
$ws.Range("E3").Value = '  +5.75%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.85'
$ws.Range("E5").Value = '  +1.26%  '
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.657'
$ws.Range("E6").Value = '  +2.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '70.84'
$ws.Range("E7").Value = '  +12.38%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.463'
$ws.Range("E9").Value = '  +3.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0971'
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.17'
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("B12").Value = 'Avalanche'
$ws.Range("C12").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.52'
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.724.92'
$ws.Range("E13").Value = '  +5.52%  '
$ws.Range("E14").Value = '  -0.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.77'
$ws.Range("E15").Value = '  +1.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.25'
$ws.Range("E16").Value = '  +2.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.857'
$ws.Range("E17").Value = '  +3.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.376.35'
$ws.Range("E18").Value = '  +5.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.433.50'
$ws.Range("E19").Value = '  -0.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0989'
$ws.Range("E20").Value = '  +1.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.35'
$ws.Range("E21").Value = '  +4.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.22'
$ws.Range("E22").Value = '  +2.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '251.67'
$ws.Range("E23").Value = '  +1.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.96'
$ws.Range("E24").Value = '  +18.85%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.46'
$ws.Range("E26").Value = '  +1.87%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.03'
$ws.Range("E27").Value = '  +9.78%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.05'
$ws.Range("E28").Value = '  +2.61%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.24'
$ws.Range("E29").Value = '  -2.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.17'
$ws.Range("E30").Value = '  +1.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.55'
$ws.Range("E31").Value = '  +9.60%  '
$ws.Range("E32").Value = '  -8.71%  '
$ws.Range("E33").Value = '  +2.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.99'
$ws.Range("E34").Value = '  +4.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0693'
$ws.Range("E35").Value = '  +1.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.10'
$ws.Range("E36").Value = '  +3.31%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.46'
$ws.Range("E37").Value = '  +8.16%  '
$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.59'
$ws.Range("E38").Value = '  +3.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.64'
$ws.Range("E39").Value = '  +0.28%  '
$ws.Range("E40").Value = '  +0.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.96'
$ws.Range("E41").Value = '  +4.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '18.52'
$ws.Range("E43").Value = '  +8.95%  '
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.52'
$ws.Range("E45").Value = '  +4.65%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '99.77'
$ws.Range("E46").Value = '  +2.47%  '
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.23'
$ws.Range("E47").Value = '  +2.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0948'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.452.56'
$ws.Range("E49").Value = '  +0.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.599.19'
$ws.Range("E50").Value = '  +5.82%  '
$ws.Range("E51").Value = '  -0.73%  '
